$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -6
$ws.Range("F3").Value = -9
$ws.Range("F5").Value = -8
$ws.Range("F7").Value = 4
